$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple D/E price+volume updates ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.996.00"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.868.38"
$ws.Range("E3").Value = "  -3.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.72"
$ws.Range("E5").Value = "  -0.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5067"
$ws.Range("E7").Value = "  -2.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3921"
$ws.Range("E8").Value = "  -1.89%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08176"
$ws.Range("E9").Value = "  -3.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.98"
$ws.Range("E10").Value = "  -2.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.57"
$ws.Range("E12").Value = "  +5.98%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.863.58"
$ws.Range("E13").Value = "  -3.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.252"
$ws.Range("E14").Value = "  -1.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.139"
$ws.Range("E15").Value = "  -3.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.97"
$ws.Range("E17").Value = "  -2.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001077"
$ws.Range("E18").Value = "  -3.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06336"
$ws.Range("E19").Value = "  -6.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.82"
$ws.Range("E20").Value = "  -1.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "29.989.80"
$ws.Range("E22").Value = "  -0.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.788"
$ws.Range("E23").Value = "  -4.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.05"
$ws.Range("E24").Value = "  -1.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.210"
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.085.82"
$ws.Range("E26").Value = "  -2.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.28"
$ws.Range("E27").Value = "  +0.79%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.84"
$ws.Range("E28").Value = "  -0.99%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.237"
$ws.Range("E29").Value = "  -9.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.67"
$ws.Range("E30").Value = "  -2.38%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1032"
$ws.Range("E31").Value = "  -2.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.039"
$ws.Range("E32").Value = "  -3.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.863"
$ws.Range("E33").Value = "  -4.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.739"
$ws.Range("E34").Value = "  +2.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02425"
$ws.Range("E35").Value = "  -3.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.204"
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06330"
$ws.Range("E37").Value = "  -4.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2137"
$ws.Range("E38").Value = "  -3.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6265"
$ws.Range("E41").Value = "  -4.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.209"
$ws.Range("E42").Value = "  -2.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.20"
$ws.Range("E43").Value = "  -2.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5882"
$ws.Range("E45").Value = "  -4.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.75"
$ws.Range("E46").Value = "  -3.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.987"
$ws.Range("E48").Value = "  -3.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.123"
$ws.Range("E51").Value = "  -2.51%  "

# --- E-only volume updates ---
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("E11").Value = "  -3.07%  "
$ws.Range("E16").Value = "  +0.26%  "
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("E47").Value = "  -2.25%  "

# --- Row 39/40 swap: ARBITRUM <-> FraxShare ---
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.517"
$ws.Range("E39").Value = "  -5.80%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.166"
$ws.Range("E40").Value = "  -6.46%  "

# --- Row 49/50 swap: EOS <-> Quant ---
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "121.82"
$ws.Range("E49").Value = "  -2.98%  "
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.205"
$ws.Range("E50").Value = "  -3.11%  "
